$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab Neo4j query in cell B2: add age/weight aliases to the
# WITH clause, convert Age and Weight (kg) to integers when they are whole
# numbers via apoc.convert.toInteger, and append a new `Cohort` column
# sourced from co.cohort_description.
$query = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
    WHERE f.file_format IN ["tif"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS Cohort
'@

$ws.Range("B2").Value2 = $query

# Row-height tweaks that accompanied the query edit (re-wrapped text needed
# slightly more vertical room in rows 2-4).
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 261

Write-Host "Updated CasesTab query and row heights."
